# Update Betfair Back/Lay odds values per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.1
$ws.Range("H2").Value = 2.74
$ws.Range("I2").Value = 3.45
$ws.Range("H5").Value = 2.06
$ws.Range("J5").Value = 3.55
$ws.Range("P5").Value = 1.98
$ws.Range("F6").Value = 1.84
$ws.Range("H6").Value = 4.2
$ws.Range("AI8").Value = 980
$ws.Range("I10").Value = 2.7
$ws.Range("J10").Value = 3.65
$ws.Range("F13").Value = 1.31
$ws.Range("X13").Value = 26
$ws.Range("H15").Value = 1.81
$ws.Range("Q16").Value = 1.77
$ws.Range("F17").Value = 1.47
$ws.Range("G17").Value = 1.52
$ws.Range("X17").Value = 25
$ws.Range("AB17").Value = 11
$ws.Range("AF17").Value = 12.5
$ws.Range("AL17").Value = 36
$ws.Range("AN17").Value = 6.8
$ws.Range("Q18").Value = 1.78
$ws.Range("F19").Value = 1.65
$ws.Range("G19").Value = 1.68
$ws.Range("H19").Value = 5.6
$ws.Range("I19").Value = 6
$ws.Range("Q19").Value = 1.84
$ws.Range("U19").Value = 2.02
$ws.Range("Z19").Value = 110
$ws.Range("AD19").Value = 23
$ws.Range("AJ19").Value = 16
$ws.Range("AL19").Value = 38
$ws.Range("AN19").Value = 9.199999999999999
$ws.Range("F20").Value = 2.36
$ws.Range("G20").Value = 2.42
$ws.Range("F21").Value = 1.38
$ws.Range("G21").Value = 1.39
$ws.Range("H21").Value = 10
$ws.Range("I21").Value = 10.5
$ws.Range("X21").Value = 26
$ws.Range("AD21").Value = 75
$ws.Range("AN21").Value = 5
$ws.Range("F22").Value = 1.65
$ws.Range("H22").Value = 5.6
$ws.Range("I22").Value = 6.4
$ws.Range("Q22").Value = 1.59
$ws.Range("G23").Value = 1.81
$ws.Range("H27").Value = 3.25
$ws.Range("I27").Value = 3.8
$ws.Range("K27").Value = 3.45
$ws.Range("G31").Value = 2.06
$ws.Range("I31").Value = 5.6
$ws.Range("Q31").Value = 1.83
$ws.Range("S31").Value = 2.94
$ws.Range("W31").Value = 1.94
